$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate Lithuanian category descriptions (column B) into English.
$ws.Range("B3").Value = "All household consumption expenses (monthly) "
$ws.Range("B4").Value = "Food and non-alcoholic beverages "
$ws.Range("B15").Value = "Non-alcoholic beverages"
$ws.Range("B25").Value = "Alcoholic beverages, tobacco, and drugs"
$ws.Range("B37").Value = "Clothing and footwear"
$ws.Range("B46").Value = "Housing, water, electricity, gas, and other fuels"
$ws.Range("B62").Value = "Furnishings, household equipment, and routine home maintenance"
$ws.Range("B82").Value = "Health"
$ws.Range("B98").Value = "Transport"
$ws.Range("B119").Value = "Information and communication"
$ws.Range("B136").Value = "Recreation, sports, and culture"
$ws.Range("B169").Value = "Education services"
$ws.Range("B180").Value = "Restaurants and accommodation services"
$ws.Range("B186").Value = "Insurance and financial services"
$ws.Range("B195").Value = "Personal care, social protection, and miscellaneous goods and services"

# Row heights reflow slightly after the re-wrap of the new English text.
$ws.Rows.Item(3).RowHeight = 46
$ws.Rows.Item(4).RowHeight = 35
$ws.Rows.Item(15).RowHeight = 24
$ws.Rows.Item(37).RowHeight = 24
$ws.Rows.Item(46).RowHeight = 57
$ws.Rows.Item(82).RowHeight = 13
$ws.Rows.Item(98).RowHeight = 13
$ws.Rows.Item(119).RowHeight = 46
$ws.Rows.Item(169).RowHeight = 24
$ws.Rows.Item(180).RowHeight = 46

# Selection moved as part of the editing session.
$ws.Range("F9").Select()
